$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'4.17%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'35.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'15.44%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.100"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'2.87%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07858"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'5.49%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'2.281"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'1.35%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'8.070"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'4.27%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'4.018"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'6.80%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.9278"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'0.87%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.09953"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'6.67%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.1820"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'5.42%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.08707"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'4.57%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.03385"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'5.37%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.09914"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.18%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001474"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.71%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.005675"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-0.96%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.485"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.096"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.63%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'3.02%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.1321"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'1.43%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'4.538"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'8.63%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.2235"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'5.39%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04663"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'3.30%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001239"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'1.72%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004488"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'5.34%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.0002699"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-20.33%"
$ws.Range("E27").ClearFormats()
$ws.Range("D39").Value = "'0.01757"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'8.44%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.04699"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'2.73%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007821"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'5.59%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'4.32%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.008792"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-10.49%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.002289"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'3.22%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.009197"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'5.58%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006056"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-0.61%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'0.15%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'5.793"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'118.21%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.002689"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'34.71%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.15%"
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'0.15%"
$ws.Range("E51").ClearFormats()
